$d = $word.ActiveDocument

# Locate the paragraph "should include default values and explanation in
# documentation. " - the new bullet is inserted directly after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "should include default values and explanation in documentation. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $rng.Paragraphs(1)
$anchorEnd = $anchorPara.Range.End

# Insert a new (initially empty) paragraph right after the anchor. It
# inherits the ListParagraph / numId=5 bullet formatting automatically,
# but we overwrite its content below with an explicit OOXML fragment so we
# get the exact run/proofErr layout we need.
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Range($anchorEnd, $anchorEnd).Paragraphs(1)

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:pPr>" +
        "<w:pStyle w:val='ListParagraph'/>" +
        "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='5'/></w:numPr>" +
        "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>" +
        "<w:t>How to capture content in ${openQuote}def __</w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r>" +
        "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>" +
        "<w:t>init</w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r>" +
        "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>" +
        "<w:t>__${closeQuote} function, which is currently ignored by Sphinx.</w:t>" +
    "</w:r>" +
    "</w:p>"

$newPara.Range.InsertXML($xml) | Out-Null
